# Auto-generated edit script applying cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells so numeric-looking strings are preserved verbatim
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = '26.948.05'
$ws.Range("E2").Value = '  +1.92%  '

$ws.Range("D3").Value = '1.813.02'
$ws.Range("E3").Value = '  +3.28%  '

$ws.Range("E4").Value = '  +1.04%  '

$ws.Range("D5").Value = '312.43'
$ws.Range("E5").Value = '  +2.20%  '

$ws.Range("E6").Value = '  +0.12%  '

$ws.Range("D7").Value = '0.4292'
$ws.Range("E7").Value = '  -0.97%  '

$ws.Range("E8").Value = '  +2.55%  '

$ws.Range("D9").Value = '0.07243'
$ws.Range("E9").Value = '  +1.99%  '

$ws.Range("D10").Value = '0.8632'
$ws.Range("E10").Value = '  +3.97%  '

$ws.Range("B11").Value = 'WrappedEther'
$ws.Range("C11").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D11").Value = '2.037.25'
$ws.Range("E11").Value = '  +15.69%  '

$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").Value = '21.23'
$ws.Range("E12").Value = '  +5.99%  '

$ws.Range("D13").Value = '6.640'
$ws.Range("E13").Value = '  +5.66%  '

$ws.Range("D14").Value = '5.391'
$ws.Range("E14").Value = '  +3.98%  '

$ws.Range("D15").Value = '0.06917'
$ws.Range("E15").Value = '  +2.06%  '

$ws.Range("D16").Value = '80.79'
$ws.Range("E16").Value = '  +2.69%  '

$ws.Range("E17").Value = '  +0.34%  '

$ws.Range("D18").Value = '0.000008929'
$ws.Range("E18").Value = '  +3.50%  '

$ws.Range("D19").Value = '1.006'
$ws.Range("E19").Value = '  +0.63%  '

$ws.Range("D20").Value = '15.18'
$ws.Range("E20").Value = '  +2.29%  '

$ws.Range("D21").Value = '26.998.86'
$ws.Range("E21").Value = '  +1.57%  '

$ws.Range("D22").Value = '5.207'
$ws.Range("E22").Value = '  +4.77%  '

$ws.Range("D23").Value = '10.96'
$ws.Range("E23").Value = '  -0.09%  '

$ws.Range("D24").Value = '2.249.42'
$ws.Range("E24").Value = '  +13.61%  '

$ws.Range("D25").Value = '154.03'
$ws.Range("E25").Value = '  +1.30%  '

$ws.Range("D26").Value = '1.886'
$ws.Range("E26").Value = '  +0.33%  '

$ws.Range("D27").Value = '18.34'
$ws.Range("E27").Value = '  +1.86%  '

$ws.Range("D28").Value = '5.238'
$ws.Range("E28").Value = '  +4.64%  '

$ws.Range("D29").Value = '1.900'
$ws.Range("E29").Value = '  +16.91%  '

$ws.Range("D30").Value = '115.15'
$ws.Range("E30").Value = '  +1.59%  '

$ws.Range("D31").Value = '0.08942'
$ws.Range("E31").Value = '  +0.36%  '

$ws.Range("D32").Value = '0.7425'
$ws.Range("E32").Value = '  +4.62%  '

$ws.Range("D33").Value = '1.158'
$ws.Range("E33").Value = '  +7.73%  '

$ws.Range("D34").Value = '4.423'
$ws.Range("E34").Value = '  +3.36%  '

$ws.Range("D35").Value = '2.808'
$ws.Range("E35").Value = '  +1.57%  '

$ws.Range("E36").Value = '  +0.71%  '

$ws.Range("E37").Value = '  +5.74%  '

$ws.Range("D38").Value = '0.05221'
$ws.Range("E38").Value = '  +3.37%  '

$ws.Range("D39").Value = '0.01921'
$ws.Range("E39").Value = '  +2.35%  '

$ws.Range("D40").Value = '0.5094'
$ws.Range("E40").Value = '  +4.30%  '

$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = '2.742'
$ws.Range("E41").Value = '  +7.62%  '

$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '0.1646'
$ws.Range("E42").Value = '  +3.72%  '

$ws.Range("D43").Value = '6.449'
$ws.Range("E43").Value = '  +5.20%  '

$ws.Range("D44").Value = '8.278'
$ws.Range("E44").Value = '  +5.63%  '

$ws.Range("D45").Value = '107.11'
$ws.Range("E45").Value = '  +3.00%  '

$ws.Range("D46").Value = '10.38'
$ws.Range("E46").Value = '  +4.18%  '

$ws.Range("D47").Value = '1.007'
$ws.Range("E47").Value = '  +1.17%  '

$ws.Range("E48").Value = '  +6.07%  '

$ws.Range("D49").Value = '0.4578'
$ws.Range("E49").Value = '  +2.93%  '

$ws.Range("D50").Value = '0.06283'
$ws.Range("E50").Value = '  +1.89%  '

$ws.Range("D51").Value = '1.800'
$ws.Range("E51").Value = '  +6.83%  '
